$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "29.473.36"
$ws.Range("E2").Value = "  +0.16%  "
$ws.Range("D3").Value = "1.920.07"
$ws.Range("E3").Value = "  +1.01%  "
$ws.Range("D4").NumberFormat = "@"
$ws.Range("D4").Value = "1.009"
$ws.Range("E4").Value = "  +0.58%  "
$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = "325.21"
$ws.Range("E5").Value = "  -0.01%  "
$ws.Range("E6").Value = "  +0.58%  "
$ws.Range("D7").NumberFormat = "@"
$ws.Range("D7").Value = "0.4836"
$ws.Range("E7").Value = "  +0.70%  "
$ws.Range("D8").NumberFormat = "@"
$ws.Range("D8").Value = "0.4077"
$ws.Range("E8").Value = "  +0.23%  "
$ws.Range("D9").NumberFormat = "@"
$ws.Range("D9").Value = "0.08243"
$ws.Range("E9").Value = "  +2.21%  "
$ws.Range("E10").Value = "  +1.20%  "
$ws.Range("D11").NumberFormat = "@"
$ws.Range("D11").Value = "23.26"
$ws.Range("E11").Value = "  -0.46%  "
$ws.Range("D12").Value = "1.941.82"
$ws.Range("E12").Value = "  +5.17%  "
$ws.Range("D13").NumberFormat = "@"
$ws.Range("D13").Value = "6.076"
$ws.Range("E13").Value = "  +2.17%  "
$ws.Range("D14").NumberFormat = "@"
$ws.Range("D14").Value = "7.268"
$ws.Range("E14").Value = "  +2.90%  "
$ws.Range("D15").NumberFormat = "@"
$ws.Range("D15").Value = "91.84"
$ws.Range("E15").Value = "  +2.25%  "
$ws.Range("D16").NumberFormat = "@"
$ws.Range("D16").Value = "0.06892"
$ws.Range("E16").Value = "  +3.16%  "
$ws.Range("E18").Value = "  +1.04%  "
$ws.Range("D19").NumberFormat = "@"
$ws.Range("D19").Value = "17.64"
$ws.Range("E19").Value = "  +0.23%  "
$ws.Range("E20").Value = "  +0.66%  "
$ws.Range("D21").Value = "29.467.13"
$ws.Range("E21").Value = "  +0.15%  "
$ws.Range("D22").NumberFormat = "@"
$ws.Range("D22").Value = "5.672"
$ws.Range("E22").Value = "  +2.53%  "
$ws.Range("D23").NumberFormat = "@"
$ws.Range("D23").Value = "11.75"
$ws.Range("E23").Value = "  -0.14%  "
$ws.Range("D24").NumberFormat = "@"
$ws.Range("D24").Value = "2.181"
$ws.Range("E24").Value = "  +1.05%  "
$ws.Range("D25").Value = "2.147.21"
$ws.Range("E25").Value = "  +2.79%  "
$ws.Range("D26").NumberFormat = "@"
$ws.Range("D26").Value = "6.671"
$ws.Range("E26").Value = "  +9.91%  "
$ws.Range("D27").NumberFormat = "@"
$ws.Range("D27").Value = "156.08"
$ws.Range("E27").Value = "  +0.87%  "
$ws.Range("D28").NumberFormat = "@"
$ws.Range("D28").Value = "20.05"
$ws.Range("E28").Value = "  +1.26%  "
$ws.Range("E29").Value = "  +1.18%  "
$ws.Range("D30").NumberFormat = "@"
$ws.Range("D30").Value = "120.95"
$ws.Range("E30").Value = "  +2.20%  "
$ws.Range("D31").NumberFormat = "@"
$ws.Range("D31").Value = "1.019"
$ws.Range("E31").Value = "  -1.25%  "
$ws.Range("D32").NumberFormat = "@"
$ws.Range("D32").Value = "0.09615"
$ws.Range("E32").Value = "  +1.23%  "
$ws.Range("D33").NumberFormat = "@"
$ws.Range("D33").Value = "5.661"
$ws.Range("E33").Value = "  +4.68%  "
$ws.Range("D34").NumberFormat = "@"
$ws.Range("D34").Value = "3.550"
$ws.Range("E34").Value = "  +0.24%  "
$ws.Range("E35").Value = "  -1.01%  "
$ws.Range("D36").NumberFormat = "@"
$ws.Range("D36").Value = "0.02286"
$ws.Range("E36").Value = "  +1.59%  "
$ws.Range("D37").NumberFormat = "@"
$ws.Range("D37").Value = "0.06115"
$ws.Range("E37").Value = "  +0.79%  "
$ws.Range("D38").NumberFormat = "@"
$ws.Range("D38").Value = "1.185"
$ws.Range("E38").Value = "  +0.80%  "
$ws.Range("D39").NumberFormat = "@"
$ws.Range("D39").Value = "8.068"
$ws.Range("E39").Value = "  +2.53%  "
$ws.Range("D40").NumberFormat = "@"
$ws.Range("D40").Value = "0.5984"
$ws.Range("E40").Value = "  +1.87%  "
$ws.Range("D41").NumberFormat = "@"
$ws.Range("D41").Value = "10.85"
$ws.Range("E41").Value = "  +6.32%  "
$ws.Range("E42").Value = "  +0.32%  "
$ws.Range("D43").NumberFormat = "@"
$ws.Range("D43").Value = "1.282"
$ws.Range("E43").Value = "  +0.29%  "
$ws.Range("D44").NumberFormat = "@"
$ws.Range("D44").Value = "2.391"
$ws.Range("E44").Value = "  -0.64%  "
$ws.Range("B45").Value = "EnergySwap"
$ws.Range("C45").Value = "https://coinranking.com/coin/SbWqqTui-+energyswap-ens"
$ws.Range("D45").NumberFormat = "@"
$ws.Range("D45").Value = "12.49"
$ws.Range("E45").Value = "  +2.25%  "
$ws.Range("B46").Value = "Cronos"
$ws.Range("C46").Value = "https://coinranking.com/coin/65PHZTpmE55b+cronos-cro"
$ws.Range("D46").NumberFormat = "@"
$ws.Range("D46").Value = "0.07606"
$ws.Range("E46").Value = "  -2.16%  "
$ws.Range("D47").NumberFormat = "@"
$ws.Range("D47").Value = "0.5607"
$ws.Range("E47").Value = "  +1.54%  "
$ws.Range("D48").NumberFormat = "@"
$ws.Range("D48").Value = "1.955"
$ws.Range("E48").Value = "  +1.82%  "
$ws.Range("D49").NumberFormat = "@"
$ws.Range("D49").Value = "118.42"
$ws.Range("E49").Value = "  +4.38%  "
$ws.Range("E50").Value = "  +4.03%  "
$ws.Range("D51").NumberFormat = "@"
$ws.Range("D51").Value = "72.39"
$ws.Range("E51").Value = "  +0.39%  "
